# "Generate Report for Handoff"
#
# The localization status report is being regenerated: the overall
# Status changes from "Handed back: in sync with en-US" to
# "Ready for handoff", and the associated timestamps are refreshed to the
# new generation time. Because the new status text is shorter, the
# Status columns that were previously auto-sized to fit the old text are
# re-sized (narrower) to fit the new text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps ---
# Latest HO Xliff Generate Date (Overview!G2) and Latest Handoff Datetime
# for de-de (de-de!H2) share the same underlying value.
$ws1.Range("G2").Value = "2016-09-01 15:24:46"
$ws3.Range("H2").Value = "2016-09-01 15:24:46"

# Latest Handoff Datetime for zh-cn (zh-cn!H2)
$ws2.Range("H2").Value = "2016-09-01 15:24:41"

# --- Column widths: narrower now that the status text is shorter ---
$ws1.Columns.Item(5).ColumnWidth = 16.3
$ws1.Columns.Item(6).ColumnWidth = 16.3
$ws2.Columns.Item(3).ColumnWidth = 16.3
$ws3.Columns.Item(3).ColumnWidth = 16.3
